# Task 50 - страница ошибки + settings : Production
#
# Adds a new "task 50" entry (error page / Error.cshtml) to both the
# "план" sheet (row 52) and the "Что нужно делать" sheet (row 12).

$wb = $excel.ActiveWorkbook

$wsPlan = $wb.Worksheets.Item("план")
$wsTodo = $wb.Worksheets.Item("Что нужно делать")

# --- Sheet "план": new row 52 -------------------------------------------
$wsPlan.Range("A52").Value = 50
$wsPlan.Range("B52").Value = "Страница ошибки"
$wsPlan.Range("G52").Value = "Багиров"
$wsPlan.Range("K52").Value = "Error.cshtml"
$wsPlan.Rows.Item(52).RowHeight = 22.5

# --- Sheet "Что нужно делать": new row 12 -------------------------------
$wsTodo.Range("A12").Value = 10
$wsTodo.Range("B12").Value = "Страница ошибки"
$wsTodo.Range("D12").Value = "Багиров"
$wsTodo.Range("E12").Value = "task 50"
$wsTodo.Range("F12").Value = "выполнено"

# --- Selections: update per-sheet active cell without stealing the -----
# --- workbook's active tab away from "Что нужно делать" ----------------
$wsPlan.Range("B52").Select()
$wsTodo.Select()
$wsTodo.Range("F12").Select()
